# Adds two new columns, I ("I0") and J ("IF"), to the sheet, with header
# cells styled like the existing headers (bold, centered, bordered) and
# numeric data for rows 2-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy the formatting of the existing "IP" header cell (H1) onto the new
# header cells so they pick up the same bold/centered/bordered style
# already registered in the workbook, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data ----------------------------------------------------------------
$iVals = @(6,2,9,5,6,7,8,8,7,8,7,10,8,8,9,8,8,8,11,7,8,8,8,8,7,7,7,7,6,6,6,2,10,7,8,6,7,8,9,8,2,8,8,4,8,7,8,6,5,4,4,6,11,8,11,6,7,8,8,8,8,8,9,9,4,4,3,4,5)
$jVals = @(6,3,9,5,7,7,8,8,7,8,8,10,8,8,9,8,8,8,11,7,9,8,8,8,7,7,7,7,6,6,6,3,10,7,8,6,7,8,9,8,3,8,8,5,8,7,8,6,6,5,4,6,11,8,11,6,7,8,8,8,8,8,9,9,4,4,3,4,5)

for ($idx = 0; $idx -lt $iVals.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
